# expansão das análises automáticas
# Adds three new computed columns (L: apoio_medio, M: contribuicoes,
# N: media_contribuicoes) to the summary table, extending the sheet
# from A1:K16 to A1:N16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new column headers (row 1), matching the existing header formatting ---
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# --- new column data (rows 2-16) ---
$data = @(
    @(94.08517009767365, 7547, 260.2413793103448),
    @(129.5537663642677, 32860, 395.9036144578313),
    @(79.80942858649695, 48629, 347.35),
    @(87.19342470373856, 174471, 302.9010416666667),
    @(50.12134015913439, 46, 23),
    @(95.44644410600942, 15501, 224.6521739130435),
    @(96.51058441972074, 95943, 218.0522727272727),
    @(66.65034280439198, 17194, 97.69318181818181),
    @(86.43541554443971, 74806, 108.2575976845152),
    @(55.58374799260083, 202, 28.85714285714286),
    @(30.9975069667077, 37, 5.285714285714286),
    @(40.89129143626957, 25, 12.5),
    @(18.08263434560471, 307, 17.05555555555556),
    @(15.74416694302886, 591, 23.64),
    @(20.96281755102498, 1248, 12.48)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 12).Value = $vals[0]
    $ws.Cells.Item($row, 13).Value = $vals[1]
    $ws.Cells.Item($row, 14).Value = $vals[2]
}
